# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") held stale "Strike#" counts; this rewrites each row's K
# value with the freshly computed count from the regenerated save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 4
    9  = 2
    10 = 2
    11 = 3
    12 = 3
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    19 = 2
    20 = 1
    21 = 5
    22 = 4
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 2
    31 = 3
    32 = 5
    33 = 6
    34 = 2
    35 = 1
    36 = 1
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
